# Auto-generated script to apply market-data refresh values
# to the Leve profit tables across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 244.85715
$ws.Range("I2").Value = 130
$ws.Range("J2").Value = 331
$ws.Range("K2").Value = 130
$ws.Range("L2").Value = 331
$ws.Range("M2").Value = -17
$ws.Range("N2").Value = -557
$ws.Range("H88").Value = 963.2857
$ws.Range("I88").Value = 1294.6666
$ws.Range("J88").Value = 714.75
$ws.Range("K88").Value = 1294.6666
$ws.Range("L88").Value = 714.75
$ws.Range("M88").Value = -888.6666
$ws.Range("N88").Value = -1526.75
$ws.Range("H91").Value = 963.2857
$ws.Range("I91").Value = 1294.6666
$ws.Range("J91").Value = 714.75
$ws.Range("K91").Value = 1294.6666
$ws.Range("L91").Value = 714.75
$ws.Range("M91").Value = 109.3334
$ws.Range("N91").Value = -3522.75
$ws.Range("H98").Value = 998.8
$ws.Range("I98").Value = 655.5714
$ws.Range("J98").Value = 1799.6666
$ws.Range("K98").Value = 655.5714
$ws.Range("L98").Value = 1799.6666
$ws.Range("M98").Value = 842.4286
$ws.Range("N98").Value = -4795.6666
$ws.Range("H122").Value = 998.8
$ws.Range("I122").Value = 655.5714
$ws.Range("J122").Value = 1799.6666
$ws.Range("K122").Value = 1966.7142
$ws.Range("L122").Value = 5398.9998
$ws.Range("M122").Value = 483.2857999999999
$ws.Range("N122").Value = -10298.9998
$ws.Range("H131").Value = 2176.4
$ws.Range("I131").Value = 1595.5
$ws.Range("J131").Value = 4500
$ws.Range("K131").Value = 4786.5
$ws.Range("L131").Value = 13500
$ws.Range("M131").Value = 253.5
$ws.Range("N131").Value = -23580
$ws.Range("H137").Value = 1639.7858
$ws.Range("J137").Value = 2959.6
$ws.Range("L137").Value = 8878.799999999999
$ws.Range("N137").Value = -13978.8
$ws.Range("H138").Value = 2354.9333
$ws.Range("I138").Value = 1290
$ws.Range("J138").Value = 5283.5
$ws.Range("K138").Value = 3870
$ws.Range("L138").Value = 15850.5
$ws.Range("M138").Value = 1270
$ws.Range("N138").Value = -26130.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 6026
$ws.Range("I33").Value = 6026
$ws.Range("K33").Value = 6026
$ws.Range("M33").Value = -5697
$ws.Range("H45").Value = 2939.4443
$ws.Range("I45").Value = 1976.4
$ws.Range("K45").Value = 1976.4
$ws.Range("M45").Value = -1599.4
$ws.Range("H61").Value = 1660.125
$ws.Range("I61").Value = 1504.1333
$ws.Range("K61").Value = 1504.1333
$ws.Range("M61").Value = -1292.1333
$ws.Range("H74").Value = 4399.5
$ws.Range("I74").Value = 3914.6
$ws.Range("K74").Value = 3914.6
$ws.Range("M74").Value = -3040.6
$ws.Range("H77").Value = 4399.5
$ws.Range("I77").Value = 3914.6
$ws.Range("K77").Value = 19573
$ws.Range("M77").Value = -15205
$ws.Range("H136").Value = 1660.125
$ws.Range("I136").Value = 1504.1333
$ws.Range("K136").Value = 4512.3999
$ws.Range("M136").Value = -1962.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5248.625
$ws.Range("I86").Value = 2250
$ws.Range("J86").Value = 6248.1665
$ws.Range("K86").Value = 2250
$ws.Range("L86").Value = 6248.1665
$ws.Range("M86").Value = -1127
$ws.Range("N86").Value = -8494.166499999999
$ws.Range("H89").Value = 5248.625
$ws.Range("I89").Value = 2250
$ws.Range("J89").Value = 6248.1665
$ws.Range("K89").Value = 11250
$ws.Range("L89").Value = 31240.8325
$ws.Range("M89").Value = -5634
$ws.Range("N89").Value = -42472.8325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1581.1666
$ws.Range("I16").Value = 2196.3333
$ws.Range("J16").Value = 966
$ws.Range("K16").Value = 2196.3333
$ws.Range("L16").Value = 966
$ws.Range("M16").Value = -1909.3333
$ws.Range("N16").Value = -1540
$ws.Range("H58").Value = 3368.6875
$ws.Range("I58").Value = 1538.3636
$ws.Range("J58").Value = 7395.4
$ws.Range("K58").Value = 1538.3636
$ws.Range("L58").Value = 7395.4
$ws.Range("M58").Value = -1335.3636
$ws.Range("N58").Value = -7801.4
$ws.Range("H74").Value = 79382.5
$ws.Range("J74").Value = 79382.5
$ws.Range("L74").Value = 79382.5
$ws.Range("N74").Value = -81130.5
$ws.Range("H77").Value = 79382.5
$ws.Range("J77").Value = 79382.5
$ws.Range("L77").Value = 238147.5
$ws.Range("N77").Value = -246883.5
$ws.Range("H99").Value = 3360.8572
$ws.Range("J99").Value = 3105.4
$ws.Range("L99").Value = 3105.4
$ws.Range("N99").Value = -6101.4
$ws.Range("H113").Value = 1581.1666
$ws.Range("I113").Value = 2196.3333
$ws.Range("J113").Value = 966
$ws.Range("K113").Value = 2196.3333
$ws.Range("L113").Value = 966
$ws.Range("M113").Value = -26.33329999999978
$ws.Range("N113").Value = -5306
$ws.Range("H126").Value = 3360.8572
$ws.Range("J126").Value = 3105.4
$ws.Range("L126").Value = 9316.200000000001
$ws.Range("N126").Value = -14256.2
$ws.Range("H136").Value = 3368.6875
$ws.Range("I136").Value = 1538.3636
$ws.Range("J136").Value = 7395.4
$ws.Range("K136").Value = 4615.0908
$ws.Range("L136").Value = 22186.2
$ws.Range("M136").Value = -2065.0908
$ws.Range("N136").Value = -27286.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1282.5
$ws.Range("I5").Value = 997.25
$ws.Range("J5").Value = 1472.6666
$ws.Range("K5").Value = 2991.75
$ws.Range("L5").Value = 4417.9998
$ws.Range("M5").Value = -2879.75
$ws.Range("N5").Value = -4641.9998
$ws.Range("H135").Value = 1282.5
$ws.Range("I135").Value = 997.25
$ws.Range("J135").Value = 1472.6666
$ws.Range("K135").Value = 8975.25
$ws.Range("L135").Value = 13253.9994
$ws.Range("M135").Value = -6440.25
$ws.Range("N135").Value = -18323.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1219.6
$ws.Range("I80").Value = 1233
$ws.Range("K80").Value = 1233
$ws.Range("M80").Value = -235
$ws.Range("H83").Value = 1219.6
$ws.Range("I83").Value = 1233
$ws.Range("K83").Value = 6165
$ws.Range("M83").Value = -1173
$ws.Range("H132").Value = 91992.17999999999
$ws.Range("I132").Value = 101102
$ws.Range("J132").Value = 894
$ws.Range("K132").Value = 303306
$ws.Range("L132").Value = 2682
$ws.Range("M132").Value = -300776
$ws.Range("N132").Value = -7742

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1399.8
$ws.Range("I7").Value = 1399.8
$ws.Range("K7").Value = 1399.8
$ws.Range("M7").Value = -1287.8
$ws.Range("H40").Value = 6468.0835
$ws.Range("I40").Value = 5892.2383
$ws.Range("K40").Value = 5892.2383
$ws.Range("M40").Value = -5756.2383
$ws.Range("H46").Value = 5242.8096
$ws.Range("I46").Value = 4333.625
$ws.Range("K46").Value = 4333.625
$ws.Range("M46").Value = -4145.625
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H123").Value = 73619
$ws.Range("J123").Value = 73619
$ws.Range("L123").Value = 73619
$ws.Range("N123").Value = -83419
$ws.Range("H126").Value = 1399.8
$ws.Range("I126").Value = 1399.8
$ws.Range("K126").Value = 4199.4
$ws.Range("M126").Value = -1729.4
$ws.Range("H132").Value = 3413
$ws.Range("I132").Value = 3378.2144
$ws.Range("K132").Value = 10134.6432
$ws.Range("M132").Value = -7604.643199999999
$ws.Range("H136").Value = 2498.6667
$ws.Range("I136").Value = 2498.6667
$ws.Range("K136").Value = 7496.000100000001
$ws.Range("M136").Value = -4946.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1232.3334
$ws.Range("I96").Value = 1098.5
$ws.Range("K96").Value = 1098.5
$ws.Range("M96").Value = 274.5
$ws.Range("H136").Value = 2675.5
$ws.Range("I136").Value = 1819.2142
$ws.Range("K136").Value = 5457.642599999999
$ws.Range("M136").Value = -2907.642599999999

Write-Output "Applied all updates"